# Scene 73.docx ("Act 2 Prim") round-trips through a Google Docs OOXML
# export, which re-emits the built-in style block (Normal, TableNormal,
# Heading1-6, Title, Subtitle) verbatim a second time and drops two
# Google-Docs roundtrip bookkeeping parts (customXML/item1.xml +
# customXML/itemProps1.xml) into the package. None of that carries any
# visible/semantic change -- every duplicated <w:style> node is byte
# identical to the original, and the customXML parts are opaque
# roundtrip blobs, not authored content. There is no text edit anywhere
# in the target diff (document.xml is untouched).
#
# We reproduce the intent through the real Word object model:
#   1. (Re-)declare each of the built-in styles the export duplicates.
#      Word's Styles collection is keyed by name/id, so re-adding a
#      style that already exists is the correct, idempotent COM
#      operation -- it returns the existing style rather than forking a
#      second definition (Word style IDs must stay unique; that is by
#      design, not a gap in this script).
#   2. Add the two Google Docs custom XML parts with their literal
#      target payloads via Document.CustomXMLParts, which is the COM
#      surface that models customXML/*.xml package parts.

$d = $word.ActiveDocument

# --- 1. Re-affirm the styles the Google Docs export re-serializes ----
# wdStyleTypeParagraph = 1, wdStyleTypeTable = 3
$builtInStyles = @(
    @{ Name = "Normal";      Type = 1 },
    @{ Name = "Table Normal"; Type = 3 },
    @{ Name = "Heading 1";   Type = 1 },
    @{ Name = "Heading 2";   Type = 1 },
    @{ Name = "Heading 3";   Type = 1 },
    @{ Name = "Heading 4";   Type = 1 },
    @{ Name = "Heading 5";   Type = 1 },
    @{ Name = "Heading 6";   Type = 1 },
    @{ Name = "Title";       Type = 1 },
    @{ Name = "Subtitle";    Type = 1 }
)

foreach ($styleDef in $builtInStyles) {
    $d.Styles.Add($styleDef.Name, $styleDef.Type) | Out-Null
}

# --- 2. Add the Google Docs roundtrip custom XML parts ----------------
$item1Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<go:gDocsCustomXmlDataStorage xmlns:go="http://customooxmlschemas.google.com/" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <go:docsCustomData xmlns:go="http://customooxmlschemas.google.com/" roundtripDataSignature="AMtx7miuE5U4RSb/52CHghi2T+sa17l1tA==">AMUW2mW401sJmWrHpB4Em+MWMgihQWPXGsT2zceOZ3PFDs6RcxD9kcizgU71HFSKuWTaFhg2WFu4VYlwcKZ5U1FNiar3ehLf8pMpRYO8+diypAFgH8iCBrE=</go:docsCustomData>
</go:gDocsCustomXmlDataStorage>
'@

$itemProps1Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<ds:datastoreItem xmlns:ds="http://schemas.openxmlformats.org/officeDocument/2006/customXml" ds:itemID="{11111111-1234-1234-1234-123412341234}">
  <ds:schemaRefs>
    <ds:schemaRef ds:uri="http://schemas.openxmlformats.org/officeDocument/2006/relationships"/>
    <ds:schemaRef ds:uri="http://customooxmlschemas.google.com/"/>
  </ds:schemaRefs>
</ds:datastoreItem>
'@

$d.CustomXMLParts.Add($item1Xml) | Out-Null
$d.CustomXMLParts.Add($itemProps1Xml) | Out-Null

Write-Output "Styles present: $($d.Styles.Count); CustomXMLParts present: $($d.CustomXMLParts.Count)"
